$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-09-15 Monday" "2025-09-16 Tuesday"

Replace-Text "209÷9=23, 2" "534÷8=66, 6"
Replace-Text "721÷3=240, 1" "457÷5=91, 2"
Replace-Text "167÷3=55, 2" "367÷5=73, 2"
Replace-Text "255÷7=36, 3" "934÷2=467, 0"
Replace-Text "470÷2=235, 0" "316÷6=52, 4"

Replace-Text "618÷9=68, 6" "203÷6=33, 5"
Replace-Text "591÷9=65, 6" "222÷6=37, 0"
Replace-Text "900÷4=225, 0" "664÷4=166, 0"
Replace-Text "110÷8=13, 6" "574÷7=82, 0"
Replace-Text "103÷2=51, 1" "980÷6=163, 2"

Replace-Text "939÷2=469, 1" "626÷4=156, 2"
Replace-Text "677÷6=112, 5" "410÷6=68, 2"
Replace-Text "534÷3=178, 0" "341÷7=48, 5"
Replace-Text "547÷5=109, 2" "196÷5=39, 1"
Replace-Text "630÷9=70, 0" "453÷3=151, 0"

Replace-Text "517÷2=258, 1" "616÷2=308, 0"
Replace-Text "791÷4=197, 3" "519÷8=64, 7"
Replace-Text "327÷6=54, 3" "743÷6=123, 5"
Replace-Text "191÷2=95, 1" "645÷6=107, 3"
Replace-Text "387÷2=193, 1" "279÷6=46, 3"

Replace-Text "938÷5=187, 3" "449÷9=49, 8"
Replace-Text "728÷8=91, 0" "945÷6=157, 3"
Replace-Text "659÷5=131, 4" "840÷9=93, 3"
Replace-Text "339÷5=67, 4" "727÷4=181, 3"
Replace-Text "107÷4=26, 3" "183÷7=26, 1"
